# Updated symbol list on Wed Dec 28 21:51:11 UTC 2022 with GitHub Actions
#
# Refreshes the cryptocurrency price snapshot (column D) for a number of
# rows, and touches the "Worst in 24h" / "Best in 24h" marker text that is
# appended to the coin-rank label in column E for two rows.
#
# Price cells are text values (e.g. "243.55"), not numbers, so we prefix
# them with a leading apostrophe when assigning via .Value - this forces
# Excel to keep them as text (preserving formatting such as trailing
# zeros) instead of auto-converting the numeric-looking strings to real
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates -------------------------------------------
$ws.Range("D2").Value  = "'243.55"
$ws.Range("D3").Value  = "'23.80"
$ws.Range("D4").Value  = "'5.246"
$ws.Range("D5").Value  = "'0.05777"
$ws.Range("D6").Value  = "'6.432"
$ws.Range("D7").Value  = "'3.228"
$ws.Range("D8").Value  = "'0.8069"
$ws.Range("D9").Value  = "'0.8791"
$ws.Range("D10").Value = "'0.1380"
$ws.Range("D12").Value = "'0.03135"
$ws.Range("D13").Value = "'0.03037"
$ws.Range("D14").Value = "'0.09330"
$ws.Range("D15").Value = "'3.824"
$ws.Range("D17").Value = "'0.04712"
$ws.Range("D18").Value = "'0.0006035"
$ws.Range("D19").Value = "'0.006177"
$ws.Range("D20").Value = "'0.001263"
$ws.Range("D21").Value = "'0.004054"
$ws.Range("D22").Value = "'0.00008724"
$ws.Range("D23").Value = "'3.543"
$ws.Range("D24").Value = "'2.156"
$ws.Range("D26").Value = "'0.1319"
$ws.Range("D28").Value = "'0.0002334"
$ws.Range("D40").Value = "'0.03737"
$ws.Range("D41").Value = "'0.006277"
$ws.Range("D42").Value = "'0.1046"
$ws.Range("D43").Value = "'0.002504"
$ws.Range("D44").Value = "'0.007155"
$ws.Range("D45").Value = "'0.00005333"
$ws.Range("D48").Value = "'0.002494"

# --- Label (column E) updates: move the Worst/Best-in-24h marker --------
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
